# Update the earliest start time to 7:40 as per SPS meeting on 5/3.
#
# Column C ("Start") holds the proposed start time tier for each school
# ("7:30 a.m.", "8:30 a.m.", or "9:30 a.m."). The earliest tier is being
# moved ten minutes later, from 7:30 a.m. to 7:40 a.m., for every school
# currently on that tier.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$startCol = 3  # Column C = "Start"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $startCol)
    if ($cell.Value2 -eq "7:30 a.m.") {
        $cell.Value = "7:40 a.m."
    }
}

# Reflect that column C was the range worked on / selected during this edit.
$ws.Range("C1:C1048576").Select()
